# Remove pictures from INPN:
# - delete the hyperlink placed on C58 (the INPN photo URL shown there)
# - clear the INPN photo URL values from all the cells that held them
# - update the active selection to D28 (and drop the saved scroll position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink attached to C58 (points at an INPN photo URL).
# (Deleting through the worksheet-level Hyperlinks collection is what
# actually removes the relationship; Range.Hyperlinks.Delete() is a no-op.)
$hyperlinksToRemove = @()
foreach ($hl in $ws.Hyperlinks) {
    $hlRange = $hl.Range
    if ($hlRange.Row -eq 58 -and $hlRange.Column -eq 3) {
        $hyperlinksToRemove += $hl
    }
}
foreach ($hl in $hyperlinksToRemove) {
    $hl.Delete()
}

# Cells whose content was an INPN photo URL ("femelle"/"male" columns, C & D)
$cellsToClear = @(
    "D2",
    "C3","D3",
    "C15",
    "C18","D18",
    "C19",
    "C20","D20",
    "D21",
    "C22",
    "D23",
    "C26","D26",
    "C27",
    "C28",
    "C29","D29",
    "D30",
    "D38",
    "C39",
    "C41",
    "C42","D42",
    "C45","D45",
    "D49",
    "C56",
    "C58",
    "C63"
)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}

# Update the sheet's active selection / scroll position.
$ws.Range("D28").Select()
